$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Python")

# Row 38: new "Find First and Last Position of Element in Sorted Array" entry.
# Copy formatting from the last existing data row (37) so styles/number
# formats match the rest of the "Legend" table.
$ws.Range("A37:C37").Copy($ws.Range("A38:C38"))
$ws.Range("A38").Value = "*"
$ws.Range("B38").Value = "Find First and Last Position of Element in Sorted Array"
$ws.Range("C38").Value = 32
$ws.Rows.Item(38).RowHeight = 34

# Row 39: new " Sort Array By Parity II" entry.
$ws.Range("A37:C37").Copy($ws.Range("A39:C39"))
$ws.Range("A39").Value = "Y"
$ws.Range("B39").Value = " Sort Array By Parity II"
$ws.Range("C39").Value = 922

# Scroll the sheet down and move the active selection the way the author
# left it after appending the two rows.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 34
$ws.Range("D45").Select()
